$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.440906882286072
$ws.Range("B1").Value = 2.371410608291626
$ws.Range("C1").Value = 2.956190586090088
$ws.Range("D1").Value = 3.432315587997437
$ws.Range("E1").Value = 1.956582307815552
